# Edit script: updates Sri Lanka IPPU calibration input data
# Commit message: "updated tanzania data to be able to run"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1: elasticity_ippu_*_production_to_gdp rows (97-105, 110-111) ---
# Columns J-N (periods 0-4) receive new randomized calibration draws.
# Columns O-AS (periods 5-35) become a linear ramp from 0 (at period 5)
# up to 1 (at period 35): value = (period - 5) / 30
$elasticityRows = @{
    97 = @(7.352557421706344, -3.373132943588593, -6.488892781360207, 0.1138069291024017, 4.227970163916579)
    98 = @(-7.778556027190428, 0.902726020902136, 1.048935941166182, 1.383862547128134, -2.316354931466868)
    99 = @(4.498888620842355, 1.264092425302407, -0.0109603687084165, 0.7291214146152644, -1.261182291616252)
    100 = @(-2.881501143613087, -0.2968302677760844, 0.387496673938261, -4.423440663835002, 3.448999547701408)
    101 = @(0.2751575730778546, 22.09430693376326, -13.45110173287774, -9.545530866041863, 26.97708393610572)
    102 = @(5.763545076972254, 1.107057815380003, -3.775154178420026, 5.816744542977935, -2.77803372183387)
    103 = @(-10.22911202776086, -0.2945226864701558, 1.872952188971701, 1.559552282011421, -0.2041127431485946)
    104 = @(-0.0010643550799541, 1.3452863035827, -1.777287280633069, -1.706565729061001, 1.589989995486501)
    105 = @(0.8906164334319027, -0.2775917300133474, 0.1680711249933744, -1.107515871538016, 0.0487905994622913)
    110 = @(1.137224907363719, 0.4740441610949519, -3.751323708564448, 3.232817727912479, -1.65726284174617)
    111 = @(2.492445319429683, 0.2549311301145209, -5.270166123553977, 3.430857869259565, -0.4853753920328891)
}

foreach ($row in $elasticityRows.Keys) {
    $rowNum = [int]$row
    $vals = $elasticityRows[$row]
    # Columns J (10) through N (14): new randomized values
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Cells.Item($rowNum, 10 + $i).Value = $vals[$i]
    }
    # Columns O (period 5, col 15) through AS (period 35, col 45): linear ramp
    for ($t = 5; $t -le 35; $t++) {
        $colIndex = 10 + $t
        $ws.Cells.Item($rowNum, $colIndex).Value = ($t - 5) / 30
    }
}

# --- Section 2: prodinit_ippu_*_tonne rows (132-140, 147-148) ---
# Columns P (col 16) through AS (col 45) are flattened to equal column O
# (period 5, col 15), i.e. production held constant from period 5 onward
# instead of exploding/decaying exponentially.
$prodinitRows = @(132, 133, 134, 135, 136, 137, 138, 139, 140, 147, 148)

foreach ($row in $prodinitRows) {
    $oVal = $ws.Cells.Item($row, 15).Value2
    for ($col = 16; $col -le 45; $col++) {
        $ws.Cells.Item($row, $col).Value = $oVal
    }
}
